# --- Load History and view account info -------------------------------
# 1) Rename the "Filter" sheet to "CardHolderFilter" and fix its
#    "Departments" header to read "Department".
# 2) Add a new "LoadHistoryFilter" sheet (right after CardHolderFilter)
#    with Load-History / account lookup data.

$wb = $excel.ActiveWorkbook

# --- 1. CardHolderFilter -------------------------------------------------
$ws = $wb.Worksheets.Item("Filter")
$ws.Name = "CardHolderFilter"

# --- 2. New LoadHistoryFilter sheet, placed right after CardHolderFilter -
$newSheet = $wb.Worksheets.Add($null, $ws)
$newSheet.Name = "LoadHistoryFilter"

# Fill data in the same order the strings were first introduced
# (escalation, then Last4Digits, then Department) so the shared-string
# table comes out in the expected order.
$newSheet.Range("A1").Value = "EmployeeName"
$newSheet.Range("B1").Value = "EmployeeId"
$newSheet.Range("A2").Value = "escalation"
$newSheet.Range("B2").Value = 455545
$newSheet.Range("C2").Value = "Admin"
$newSheet.Range("D1").Value = "Last4Digits"
$newSheet.Range("D2").Value = 3802
$newSheet.Range("C1").Value = "Department"

# Match column widths (bestFit-style) seen on the other filter sheets.
$newSheet.Columns.Item(1).ColumnWidth = 12.996651785714286
$newSheet.Columns.Item(2).ColumnWidth = 9.711495535714286
$newSheet.Columns.Item(3).ColumnWidth = 10.285714285714286
$newSheet.Columns.Item(4).ColumnWidth = 9.141183035714286
$newSheet.Columns.Item(5).ColumnWidth = 10.141183035714286

$newSheet.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 130
$newSheet.Range("I16").Select() | Out-Null

# --- Update CardHolderFilter header + selection --------------------------
$ws.Range("C1").Value = "Department"
$ws.Range("D22").Select() | Out-Null

# Keep CardHolderFilter as the active/selected tab, as in the original file.
$ws.Activate() | Out-Null
